# Move the sentences from the "Sentences" sheet onto Sheet1 (replacing the
# old placeholder strings with the new "test three".."test eight" values),
# then remove the now-redundant "Sentences" sheet entirely.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$wsSentences = $wb.Worksheets.Item("Sentences")

$values = @("test three", "test four", "test five", "test six", "test seven", "test eight")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 1).Value = $values[$i]
}

$wsSentences.Delete()

$ws1.Select()
